# Fruta / hortaliza, semanal
# Insert a new week's worth of data (Primera + Segunda rows) right before the
# old row 1125, shifting all the following rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1125-1126 (old row 1125 becomes row 1127, etc.)
$ws.Range("A1125:A1126").EntireRow.Insert()

# New row 1125: Coliflor, Primera, fecha 45147
$ws.Cells.Item(1125, 1).Value = 8
$ws.Cells.Item(1125, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1125, 3).Value = "Coquimbo"
$ws.Cells.Item(1125, 4).Value = 45147
$ws.Cells.Item(1125, 5).Value = 4
$ws.Cells.Item(1125, 6).Value = 100112008
$ws.Cells.Item(1125, 7).Value = "Coliflor"
$ws.Cells.Item(1125, 8).Value = "Sin especificar"
$ws.Cells.Item(1125, 9).Value = "Primera"
$ws.Cells.Item(1125, 10).Value = 2000
$ws.Cells.Item(1125, 11).Value = 700
$ws.Cells.Item(1125, 12).Value = 800
$ws.Cells.Item(1125, 13).Value = 750
$ws.Cells.Item(1125, 14).Value = "$/unidad"
$ws.Cells.Item(1125, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1125, 16).Value = 750
$ws.Cells.Item(1125, 17).Value = 1
$ws.Cells.Item(1125, 18).Value = "Hortaliza"

# New row 1126: Coliflor, Segunda, fecha 45147
$ws.Cells.Item(1126, 1).Value = 8
$ws.Cells.Item(1126, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1126, 3).Value = "Coquimbo"
$ws.Cells.Item(1126, 4).Value = 45147
$ws.Cells.Item(1126, 5).Value = 4
$ws.Cells.Item(1126, 6).Value = 100112008
$ws.Cells.Item(1126, 7).Value = "Coliflor"
$ws.Cells.Item(1126, 8).Value = "Sin especificar"
$ws.Cells.Item(1126, 9).Value = "Segunda"
$ws.Cells.Item(1126, 10).Value = 1000
$ws.Cells.Item(1126, 11).Value = 500
$ws.Cells.Item(1126, 12).Value = 600
$ws.Cells.Item(1126, 13).Value = 550
$ws.Cells.Item(1126, 14).Value = "$/unidad"
$ws.Cells.Item(1126, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1126, 16).Value = 550
$ws.Cells.Item(1126, 17).Value = 1
$ws.Cells.Item(1126, 18).Value = "Hortaliza"
